# Auto-generated edit script: updates cryptos list (Wed Aug 30 09:29:37 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.516.33'
$ws.Range('E2').Value = '  +5.66%  '
$ws.Range('D3').Value = '1.721.62'
$ws.Range('E3').Value = '  +4.30%  '
$ws.Range('D4').Value = '''1.006'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').Value = '''225.80'
$ws.Range('E5').Value = '  +3.41%  '
$ws.Range('D6').Value = '''0.5351'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D7').Value = '''1.005'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''0.2665'
$ws.Range('E8').Value = '  +1.15%  '
$ws.Range('D9').Value = '''0.06588'
$ws.Range('E9').Value = '  +4.17%  '
$ws.Range('D10').Value = '''21.64'
$ws.Range('E10').Value = '  +6.08%  '
$ws.Range('D11').Value = '''0.07744'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').Value = '''4.618'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '1.720.31'
$ws.Range('E13').Value = '  +4.31%  '
$ws.Range('D14').Value = '1.961.62'
$ws.Range('E14').Value = '  +4.48%  '
$ws.Range('D15').Value = '''0.5837'
$ws.Range('E15').Value = '  +4.30%  '
$ws.Range('D16').Value = '0.0₅8307'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('D17').Value = '''67.88'
$ws.Range('E17').Value = '  +4.09%  '
$ws.Range('D18').Value = '27.534.25'
$ws.Range('E18').Value = '  +5.77%  '
$ws.Range('D19').Value = '''219.97'
$ws.Range('E19').Value = '  +14.79%  '
$ws.Range('D20').Value = '''1.005'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').Value = '''4.725'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('D22').Value = '''10.62'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('D23').Value = '''6.075'
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('D24').Value = '''1.006'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').Value = '''148.62'
$ws.Range('E25').Value = '  +3.57%  '
$ws.Range('D26').Value = '''1.730'
$ws.Range('E26').Value = '  +14.81%  '
$ws.Range('D27').Value = '''0.1234'
$ws.Range('E27').Value = '  +4.09%  '
$ws.Range('D28').Value = '''7.401'
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('D29').Value = '''16.57'
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').Value = '''0.05564'
$ws.Range('E30').Value = '  +3.41%  '
$ws.Range('D31').Value = '''1.303'
$ws.Range('E31').Value = '  +2.70%  '
$ws.Range('D32').Value = '''3.565'
$ws.Range('E32').Value = '  +3.34%  '
$ws.Range('D33').Value = '''3.446'
$ws.Range('E33').Value = '  +2.78%  '
$ws.Range('D34').Value = '''1.659'
$ws.Range('E34').Value = '  +6.77%  '
$ws.Range('D35').Value = '''2.858'
$ws.Range('E35').Value = '  +2.70%  '
$ws.Range('D36').Value = '''0.9660'
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('D37').Value = '''2.423'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '''0.5958'
$ws.Range('E38').Value = '  +5.77%  '
$ws.Range('D39').Value = '''0.01651'
$ws.Range('E39').Value = '  +4.73%  '
$ws.Range('D40').Value = '''5.913'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '''0.8544'
$ws.Range('E41').Value = '  +3.51%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.056.29'
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('D43').Value = '''1.005'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').Value = '''101.31'
$ws.Range('E44').Value = '  +0.43%  '
$ws.Range('D45').Value = '1.866.89'
$ws.Range('E45').Value = '  +4.48%  '
$ws.Range('D46').Value = '0.0₈111'
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('D47').Value = '''58.91'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('D48').Value = '''8.236'
$ws.Range('E48').Value = '  +4.10%  '
$ws.Range('D49').Value = '''0.4440'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').Value = '''1.005'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').Value = '''0.05253'
